# Auto-generated edit script applying the Sheets/Midgardsormr_Profits.xlsx diff
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) across the
# ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit tables per the scheduled-runner refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 4405.81
$ws.Range("I15").Value = 4405.81
$ws.Range("K15").Value = 13217.43
$ws.Range("M15").Value = -13048.43

$ws.Range("H33").Value = 937.62964
$ws.Range("I33").Value = 826.913
$ws.Range("K33").Value = 826.913
$ws.Range("M33").Value = -597.913

$ws.Range("H57").Value = 54895
$ws.Range("I57").Value = 19800
$ws.Range("J57").Value = 89990
$ws.Range("K57").Value = 59400
$ws.Range("L57").Value = 269970
$ws.Range("M57").Value = -58901
$ws.Range("N57").Value = -270968

$ws.Range("H76").Value = 4356.591
$ws.Range("I76").Value = 3626.3845
$ws.Range("J76").Value = 5411.3335
$ws.Range("K76").Value = 3626.3845
$ws.Range("L76").Value = 5411.3335
$ws.Range("M76").Value = -3311.3845
$ws.Range("N76").Value = -6041.3335

$ws.Range("H79").Value = 4356.591
$ws.Range("I79").Value = 3626.3845
$ws.Range("J79").Value = 5411.3335
$ws.Range("K79").Value = 3626.3845
$ws.Range("L79").Value = 5411.3335
$ws.Range("M79").Value = -2534.3845
$ws.Range("N79").Value = -7595.3335

$ws.Range("H137").Value = 9636.075000000001
$ws.Range("I137").Value = 22213
$ws.Range("J137").Value = 2089.92
$ws.Range("K137").Value = 66639
$ws.Range("L137").Value = 6269.76
$ws.Range("M137").Value = -64089
$ws.Range("N137").Value = -11369.76

$ws.Range("H138").Value = 3549.4722
$ws.Range("I138").Value = 2982.4614
$ws.Range("J138").Value = 3869.9565
$ws.Range("K138").Value = 8947.3842
$ws.Range("L138").Value = 11609.8695
$ws.Range("M138").Value = -3807.3842
$ws.Range("N138").Value = -21889.8695

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14658.819
$ws.Range("I32").Value = 14894.5
$ws.Range("J32").Value = 8374
$ws.Range("K32").Value = 14894.5
$ws.Range("L32").Value = 8374
$ws.Range("M32").Value = -14607.5
$ws.Range("N32").Value = -8948

$ws.Range("H45").Value = 3130
$ws.Range("I45").Value = 2329.318
$ws.Range("K45").Value = 2329.318
$ws.Range("M45").Value = -1952.318

$ws.Range("H61").Value = 3685.4285
$ws.Range("I61").Value = 3199.6
$ws.Range("K61").Value = 3199.6
$ws.Range("M61").Value = -2987.6

$ws.Range("H74").Value = 401255.88
$ws.Range("I74").Value = 429774.16
$ws.Range("K74").Value = 429774.16
$ws.Range("M74").Value = -428900.16

$ws.Range("H77").Value = 401255.88
$ws.Range("I77").Value = 429774.16
$ws.Range("K77").Value = 2148870.8
$ws.Range("M77").Value = -2144502.8

$ws.Range("H88").Value = 5464.75
$ws.Range("I88").Value = 1481.3334
$ws.Range("J88").Value = 6792.5557
$ws.Range("K88").Value = 1481.3334
$ws.Range("L88").Value = 6792.5557
$ws.Range("M88").Value = -1075.3334
$ws.Range("N88").Value = -7604.5557

$ws.Range("H91").Value = 5464.75
$ws.Range("I91").Value = 1481.3334
$ws.Range("J91").Value = 6792.5557
$ws.Range("K91").Value = 1481.3334
$ws.Range("L91").Value = 6792.5557
$ws.Range("M91").Value = -77.33339999999998
$ws.Range("N91").Value = -9600.555700000001

$ws.Range("H97").Value = 1763.5814
$ws.Range("J97").Value = 2400.9473
$ws.Range("L97").Value = 2400.9473
$ws.Range("N97").Value = -3392.9473

$ws.Range("H102").Value = 3311.279
$ws.Range("I102").Value = 2711.6365
$ws.Range("K102").Value = 2711.6365
$ws.Range("M102").Value = -1089.6365

$ws.Range("H132").Value = 2482.0356
$ws.Range("I132").Value = 1519.091
$ws.Range("K132").Value = 4557.272999999999
$ws.Range("M132").Value = -2027.272999999999

$ws.Range("H136").Value = 3685.4285
$ws.Range("I136").Value = 3199.6
$ws.Range("K136").Value = 9598.799999999999
$ws.Range("M136").Value = -7048.799999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1007.46155
$ws.Range("I94").Value = 980.2381
$ws.Range("K94").Value = 980.2381
$ws.Range("M94").Value = -529.2381

$ws.Range("H105").Value = 3221.8
$ws.Range("I105").Value = 1964.1428
$ws.Range("J105").Value = 6156.3335
$ws.Range("K105").Value = 1964.1428
$ws.Range("L105").Value = 6156.3335
$ws.Range("M105").Value = -217.1428000000001
$ws.Range("N105").Value = -9650.333500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3803.1875
$ws.Range("I62").Value = 2275.3333
$ws.Range("K62").Value = 2275.3333
$ws.Range("M62").Value = -1651.3333

$ws.Range("H65").Value = 3803.1875
$ws.Range("I65").Value = 2275.3333
$ws.Range("K65").Value = 11376.6665
$ws.Range("M65").Value = -8256.666499999999

$ws.Range("H107").Value = 647.7692
$ws.Range("I107").Value = 544.12
$ws.Range("K107").Value = 544.12
$ws.Range("M107").Value = 1375.88

$ws.Range("H122").Value = 1506.8182
$ws.Range("I122").Value = 1506.8182
$ws.Range("K122").Value = 4520.4546
$ws.Range("M122").Value = -2070.4546

$ws.Range("H132").Value = 35575.39
$ws.Range("I132").Value = 42024.535
$ws.Range("J132").Value = 3329.6667
$ws.Range("K132").Value = 126073.605
$ws.Range("L132").Value = 9989.000100000001
$ws.Range("M132").Value = -123543.605
$ws.Range("N132").Value = -15049.0001

$ws.Range("H133").Value = 147383.6
$ws.Range("J133").Value = 165442
$ws.Range("L133").Value = 165442
$ws.Range("N133").Value = -170502

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 970
$ws.Range("I117").Value = 1221.4286
$ws.Range("J117").Value = 718.5714
$ws.Range("K117").Value = 3664.2858
$ws.Range("L117").Value = 2155.7142
$ws.Range("M117").Value = -222.2857999999997
$ws.Range("N117").Value = -9039.7142

$ws.Range("H131").Value = 252597.06
$ws.Range("J131").Value = 3131.9167
$ws.Range("L131").Value = 9395.750100000001
$ws.Range("N131").Value = -19475.7501

$ws.Range("H132").Value = 1726.875
$ws.Range("I132").Value = 1866.8
$ws.Range("J132").Value = 1626.9286
$ws.Range("K132").Value = 16801.2
$ws.Range("L132").Value = 14642.3574
$ws.Range("M132").Value = -14271.2
$ws.Range("N132").Value = -19702.3574

$ws.Range("H138").Value = 15752.7
$ws.Range("I138").Value = 15646.857
$ws.Range("K138").Value = 46940.571
$ws.Range("M138").Value = -41800.571

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 29998.334
$ws.Range("J57").Value = 29998.334
$ws.Range("L57").Value = 29998.334
$ws.Range("N57").Value = -31638.334

$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H33").Value = 32250
$ws.Range("J33").Value = 32250
$ws.Range("L33").Value = 32250
$ws.Range("N33").Value = -32830

$ws.Range("H61").Value = 985.3
$ws.Range("I61").Value = 846
$ws.Range("J61").Value = 1542.5
$ws.Range("K61").Value = 846
$ws.Range("L61").Value = 1542.5
$ws.Range("M61").Value = -644
$ws.Range("N61").Value = -1946.5

$ws.Range("H113").Value = 985.3
$ws.Range("I113").Value = 846
$ws.Range("J113").Value = 1542.5
$ws.Range("K113").Value = 846
$ws.Range("L113").Value = 1542.5
$ws.Range("M113").Value = 1324
$ws.Range("N113").Value = -5882.5

$ws.Range("H122").Value = 8685.394
$ws.Range("I122").Value = 9656.35
$ws.Range("J122").Value = 7191.615
$ws.Range("K122").Value = 28969.05
$ws.Range("L122").Value = 21574.845
$ws.Range("M122").Value = -26519.05
$ws.Range("N122").Value = -26474.845

$ws.Range("H136").Value = 8329.625
$ws.Range("I136").Value = 7773
$ws.Range("J136").Value = 9999.5
$ws.Range("K136").Value = 23319
$ws.Range("L136").Value = 29998.5
$ws.Range("M136").Value = -20769
$ws.Range("N136").Value = -35098.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 37870.625
$ws.Range("I2").Value = 45050
$ws.Range("J2").Value = 16332.5
$ws.Range("K2").Value = 45050
$ws.Range("L2").Value = 16332.5
$ws.Range("M2").Value = -44938
$ws.Range("N2").Value = -16556.5

$ws.Range("H4").Value = 835013.8
$ws.Range("I4").Value = 1277.6666
$ws.Range("J4").Value = 1668750
$ws.Range("K4").Value = 1277.6666
$ws.Range("L4").Value = 1668750
$ws.Range("M4").Value = -1164.6666
$ws.Range("N4").Value = -1668976

$ws.Range("H81").Value = 3841.05
$ws.Range("J81").Value = 2516
$ws.Range("L81").Value = 5032
$ws.Range("N81").Value = -7154

$ws.Range("H84").Value = 3841.05
$ws.Range("J84").Value = 2516
$ws.Range("L84").Value = 25160
$ws.Range("N84").Value = -35768

$ws.Range("H122").Value = 58331.49
$ws.Range("I122").Value = 71372.664
$ws.Range("K122").Value = 214117.992
$ws.Range("M122").Value = -211667.992

$ws.Range("H136").Value = 73098.78
$ws.Range("I136").Value = 73098.78
$ws.Range("K136").Value = 219296.34
$ws.Range("M136").Value = -216746.34

$ws.Range("H139").Value = 132970
$ws.Range("J139").Value = 132970
$ws.Range("L139").Value = 132970
$ws.Range("N139").Value = -143250
